# Updates cryptos list (Price and Volume(1h) columns) per GitHub Actions refresh job
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.411.83"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").Value = "2.553.68"
$ws.Range("E3").Value = "  -2.25%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.81"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.27"
$ws.Range("E6").Value = "  +5.24%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").Value = "  +0.45%  "

$ws.Range("D9").Value = "2.552.22"
$ws.Range("E9").Value = "  -2.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("E10").Value = "  +1.59%  "

$ws.Range("E11").Value = "  +1.92%  "

$ws.Range("E12").Value = "  -2.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.18"
$ws.Range("E13").Value = "  -0.49%  "

$ws.Range("E14").Value = "  -0.25%  "

$ws.Range("D15").Value = "3.014.19"
$ws.Range("E15").Value = "  -2.43%  "

$ws.Range("E16").Value = "  +0.13%  "

$ws.Range("D17").Value = "67.239.12"
$ws.Range("E17").Value = "  +0.30%  "

$ws.Range("D18").Value = "2.557.03"
$ws.Range("E18").Value = "  -2.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.94"
$ws.Range("E19").Value = "  +2.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.44"
$ws.Range("E20").Value = "  -2.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "356.44"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("E22").Value = "  -0.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.69"
$ws.Range("E23").Value = "  +1.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("E24").Value = "  +6.76%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.26"
$ws.Range("E26").Value = "  +1.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.09"
$ws.Range("E27").Value = "  -3.15%  "

$ws.Range("D28").Value = "2.686.92"
$ws.Range("E28").Value = "  -2.24%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  +1.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "535.34"
$ws.Range("E31").Value = "  -1.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.28"
$ws.Range("E32").Value = "  +5.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.37"
$ws.Range("E33").Value = "  +2.40%  "

$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("E35").Value = "  +0.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("E37").Value = "  +0.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.35"
$ws.Range("E38").Value = "  -0.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.83"

$ws.Range("E40").Value = "  +1.31%  "

$ws.Range("E41").Value = "  -1.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.23"
$ws.Range("E42").Value = "  +2.26%  "

$ws.Range("E43").Value = "  +0.78%  "

$ws.Range("E44").Value = "  +6.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.75"
$ws.Range("E46").Value = "  -1.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.40"
$ws.Range("E47").Value = "  +0.34%  "

$ws.Range("D49").Value = "0.0₆0282"
$ws.Range("E49").Value = "  -4.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.73"

$ws.Range("E51").Value = "  +1.80%  "
